$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.48283851146698
$ws.Range("B1").Value = 1.265719175338745
$ws.Range("C1").Value = 5.14802360534668
$ws.Range("D1").Value = 3.51827073097229
$ws.Range("E1").Value = 0.6542774438858032
